$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (K1) onto the new
# header cells (L1:N1) so they pick up the bold/centered/bordered style.
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122) # xlPasteFormats

# New headers
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# New data values for rows 2-7
$ws.Range("L2").Value = 94.10714337626324
$ws.Range("M2").Value = 154974
$ws.Range("N2").Value = 308.0994035785288

$ws.Range("L3").Value = 87.3015904523896
$ws.Range("M3").Value = 108579
$ws.Range("N3").Value = 332.045871559633

$ws.Range("L4").Value = 86.43287869528324
$ws.Range("M4").Value = 88772
$ws.Range("N4").Value = 125.030985915493

$ws.Range("L5").Value = 93.05249603932496
$ws.Range("M5").Value = 114874
$ws.Range("N5").Value = 170.6894502228826

$ws.Range("L6").Value = 18.6886320009998
$ws.Range("M6").Value = 1789
$ws.Range("N6").Value = 14.1984126984127

$ws.Range("L7").Value = 23.27683795941807
$ws.Range("M7").Value = 419
$ws.Range("N7").Value = 16.11538461538462
